# Update cryptocurrency price (D) and 1h volume change (E) columns
# to reflect the latest scrape, matching the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.462.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.834.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5337"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.78%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3018"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06855"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7357"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.846.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07264"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.960"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.16%  "

$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007844"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.474.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.076.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.557"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.947"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.187"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.185"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("E27").Value = "  +0.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.214"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08780"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.988"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04782"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.925"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7270"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.88%  "

$ws.Range("E36").Value = "  -0.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.089"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.269"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01699"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4704"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9023"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "107.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.862"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.73%  "

$ws.Range("E44").Value = "  +0.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.333"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.912"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1227"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.16%  "

$ws.Range("E48").Value = "  -3.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05793"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8896"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.01%  "

